$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 9,20
$data[0,0] = "ECs"
$data[0,1] = "Icam5"
$data[0,2] = "Itgb2"
$data[0,3] = "ECs"
$data[0,4] = 3
$data[0,5] = 1
$data[0,6] = 1.523319666666667
$data[0,7] = 4.569959
$data[0,8] = 0.3270148067479001
$data[0,9] = 0.3270148067479002
$data[0,10] = 1
$data[0,11] = 0.3333333333333333
$data[0,12] = 0.1145113333333333
$data[0,13] = 0.343534
$data[0,14] = 0.001785365609625045
$data[0,15] = 0.001785365609625044
$data[0,16] = 0.1744373661228889
$data[0,17] = 1.569936295106
$data[0,18] = 0.0005838409898058809
$data[0,19] = 0.0005838409898058809
$data[1,0] = "ECs"
$data[1,1] = "Icam5"
$data[1,2] = "Itgb2"
$data[1,3] = "FAPs"
$data[1,4] = 3
$data[1,5] = 1
$data[1,6] = 1.523319666666667
$data[1,7] = 4.569959
$data[1,8] = 0.3270148067479001
$data[1,9] = 0.3270148067479002
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 0.467525
$data[1,13] = 1.402575
$data[1,14] = 0.007289261528465441
$data[1,15] = 0.007289261528465441
$data[1,16] = 0.7121900271583333
$data[1,17] = 6.409710244425001
$data[1,18] = 0.002383696450066029
$data[1,19] = 0.002383696450066029
$data[2,0] = "ECs"
$data[2,1] = "Icam5"
$data[2,2] = "Itgb2"
$data[2,3] = "Resolving-Mac"
$data[2,4] = 3
$data[2,5] = 1
$data[2,6] = 1.523319666666667
$data[2,7] = 4.569959
$data[2,8] = 0.3270148067479001
$data[2,9] = 0.3270148067479002
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 63.556834
$data[2,13] = 190.670502
$data[2,14] = 0.9909253728619096
$data[2,15] = 0.9909253728619095
$data[2,16] = 96.81737518326867
$data[2,17] = 871.3563766494179
$data[2,18] = 0.3240472693080282
$data[2,19] = 0.3240472693080282
$data[3,0] = "FAPs"
$data[3,1] = "Icam5"
$data[3,2] = "Itgb2"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 2.580934333333333
$data[3,7] = 7.742803
$data[3,8] = 0.5540555673983205
$data[3,9] = 0.5540555673983206
$data[3,10] = 1
$data[3,11] = 0.3333333333333333
$data[3,12] = 0.1145113333333333
$data[3,13] = 0.343534
$data[3,14] = 0.001785365609625045
$data[3,15] = 0.001785365609625044
$data[3,16] = 0.2955462317557778
$data[3,17] = 2.659916085802
$data[3,18] = 0.0009891917558542524
$data[3,19] = 0.0009891917558542526
$data[4,0] = "FAPs"
$data[4,1] = "Icam5"
$data[4,2] = "Itgb2"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 2.580934333333333
$data[4,7] = 7.742803
$data[4,8] = 0.5540555673983205
$data[4,9] = 0.5540555673983206
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 0.467525
$data[4,13] = 1.402575
$data[4,14] = 0.007289261528465441
$data[4,15] = 0.007289261528465441
$data[4,16] = 1.206651324191667
$data[4,17] = 10.859861917725
$data[4,18] = 0.004038655932068669
$data[4,19] = 0.00403865593206867
$data[5,0] = "FAPs"
$data[5,1] = "Icam5"
$data[5,2] = "Itgb2"
$data[5,3] = "Resolving-Mac"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 2.580934333333333
$data[5,7] = 7.742803
$data[5,8] = 0.5540555673983205
$data[5,9] = 0.5540555673983206
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 63.556834
$data[5,13] = 190.670502
$data[5,14] = 0.9909253728619096
$data[5,15] = 0.9909253728619095
$data[5,16] = 164.0360149885674
$data[5,17] = 1476.324134897106
$data[5,18] = 0.5490277197103977
$data[5,19] = 0.5490277197103977
$data[6,0] = "Resolving-Mac"
$data[6,1] = "Icam5"
$data[6,2] = "Itgb2"
$data[6,3] = "ECs"
$data[6,4] = 3
$data[6,5] = 1
$data[6,6] = 0.554005
$data[6,7] = 1.662015
$data[6,8] = 0.1189296258537793
$data[6,9] = 0.1189296258537793
$data[6,10] = 1
$data[6,11] = 0.3333333333333333
$data[6,12] = 0.1145113333333333
$data[6,13] = 0.343534
$data[6,14] = 0.001785365609625045
$data[6,15] = 0.001785365609625044
$data[6,16] = 0.06343985122333333
$data[6,17] = 0.5709586610099999
$data[6,18] = 0.0002123328639649111
$data[6,19] = 0.0002123328639649111
$data[7,0] = "Resolving-Mac"
$data[7,1] = "Icam5"
$data[7,2] = "Itgb2"
$data[7,3] = "FAPs"
$data[7,4] = 3
$data[7,5] = 1
$data[7,6] = 0.554005
$data[7,7] = 1.662015
$data[7,8] = 0.1189296258537793
$data[7,9] = 0.1189296258537793
$data[7,10] = 3
$data[7,11] = 1
$data[7,12] = 0.467525
$data[7,13] = 1.402575
$data[7,14] = 0.007289261528465441
$data[7,15] = 0.007289261528465441
$data[7,16] = 0.259011187625
$data[7,17] = 2.331100688625
$data[7,18] = 0.000866909146330742
$data[7,19] = 0.000866909146330742
$data[8,0] = "Resolving-Mac"
$data[8,1] = "Icam5"
$data[8,2] = "Itgb2"
$data[8,3] = "Resolving-Mac"
$data[8,4] = 3
$data[8,5] = 1
$data[8,6] = 0.554005
$data[8,7] = 1.662015
$data[8,8] = 0.1189296258537793
$data[8,9] = 0.1189296258537793
$data[8,10] = 3
$data[8,11] = 1
$data[8,12] = 63.556834
$data[8,13] = 190.670502
$data[8,14] = 0.9909253728619096
$data[8,15] = 0.9909253728619095
$data[8,16] = 35.21080382017
$data[8,17] = 316.89723438153
$data[8,18] = 0.1178503838434836
$data[8,19] = 0.1178503838434836

$ws.Range("A2:T10").Value = $data

$ws.Rows("11:13").Delete()
